# modification nom antiobiotique en francais
# - remove the obsolete "amoxicillin / G 5% / 50 mg/mL / 1 h" row
# - translate antibiotic (molecule) names in column A from English to French
# - widen columns A, B, C and E for readability
# - move the active selection to the last data row

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- 1. Delete the row that held: amoxicillin | G 5% | 50 mg/mL | 1 h | | [9] ---
$target = $null
for ($r = $ws.UsedRange.Rows.Count; $r -ge 2; $r--) {
    $a = $ws.Cells.Item($r, 1).Value2
    $b = $ws.Cells.Item($r, 2).Value2
    $c = $ws.Cells.Item($r, 3).Value2
    $d = $ws.Cells.Item($r, 4).Value2
    if ($a -eq "amoxicillin" -and $b -eq "G 5%" -and $c -eq "50 mg/mL" -and $d -eq "1 h") {
        $target = $r
        break
    }
}
if ($target) {
    $ws.Rows.Item($target).Delete()
}

# --- 2. Translate molecule names (column A) from English to French ---
$translations = @{
    "amoxicillin"                      = "amoxicilline"
    "amoxicillin + clavulanic acid"    = "amoxicilline + clavulanic acid"
    "ampicillin + sulbactam"           = "ampicilline + sulbactam"
    "cefazolin"                        = "cefazoline"
    "cefoxitin"                        = "cefoxitine"
    "clindamycin"                      = "clindamycine"
    "cloxacillin"                      = "cloxacilline"
    "colistin"                         = "colistine"
    "fosfomycin"                       = "fosfomycine"
    "oxacillin"                        = "oxacilline"
    "piperacillin"                     = "piperacilline"
    "piperacillin + tazobactam"        = "piperacilline + tazobactam"
    "teicoplanin"                      = "teicoplanine"
    "temocillin"                       = "temocilline"
    "vancomycin"                       = "vancomycine"
}

$lastRow = $ws.UsedRange.Rows.Count
for ($r = 2; $r -le $lastRow; $r++) {
    $cell = $ws.Cells.Item($r, 1)
    $val = $cell.Value2
    if ($val -ne $null -and $translations.ContainsKey($val)) {
        $cell.Value = $translations[$val]
    }
}

# --- 3. Column widths ---
$ws.Columns.Item(1).ColumnWidth = 39.1666666666667
$ws.Columns.Item(2).ColumnWidth = 17.6666666666667
$ws.Columns.Item(3).ColumnWidth = 17.6666666666667
$ws.Columns.Item(5).ColumnWidth = 29.6666666666667

# --- 4. View / selection state ---
$ws.Activate()
$lastRow = $ws.UsedRange.Rows.Count
$ws.Range("A" + $lastRow).Select()
try {
    $excel.ActiveWindow.ScrollRow = 17
} catch {
}
